$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 114: Sept 20, 2021 titration data (CRM opened 9/17/2021) ---
$ws.Range("A113").Copy()
$ws.Range("A114").PasteSpecial(-4122)   # xlPasteFormats - reuse A-column date style
$ws.Range("A114").Value = 44459
$ws.Range("B114").Value = 2238.28837133754
$ws.Range("C114").Value = 2230.52
$ws.Range("E114").Value = 183
$ws.Range("F114").Value = "CRM opened 9/17/2021"

# --- Row 115: Sept 21, 2021 titration data (CRM opened 9/17/2021) ---
$ws.Range("A113").Copy()
$ws.Range("A115").PasteSpecial(-4122)
$ws.Range("A115").Value = 44460
$ws.Range("B115").Value = 2245.6299964466002
$ws.Range("C115").Value = 2230.52
$ws.Range("E115").Value = 183
$ws.Range("F115").Value = "CRM opened 9/17/2021"

$excel.CutCopyMode = $false

# Fill in the % off formula (shared across D110:D115, continuing the existing pattern)
$ws.Range("D110:D115").Formula = "=100*(B110-C110)/C110"

# Update the frozen-pane / view state to scroll toward the new bottom rows
$ws.Range("A98").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 98
$ws.Range("A116").Select()
